$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.999.56"
$ws.Range("E2").Value = "  +2.15%  "
$ws.Range("D3").Value = "2.049.84"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'229.39"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").Value = "'0.617"
$ws.Range("E6").Value = "  +2.12%  "
$ws.Range("D7").Value = "'58.49"
$ws.Range("E7").Value = "  +5.79%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("D10").Value = "'0.0809"
$ws.Range("E10").Value = "  +2.73%  "
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("D12").Value = "2.353.24"
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("E13").Value = "  +2.44%  "
$ws.Range("E14").Value = "  +2.50%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'5.28"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.751"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").Value = "2.050.44"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").Value = "37.915.34"
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("D19").Value = "'6.24"
$ws.Range("E19").Value = "  -4.69%  "
$ws.Range("D20").Value = "'69.67"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("E21").Value = "  +2.15%  "
$ws.Range("D22").Value = "'224.66"
$ws.Range("E22").Value = "  +0.87%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'2.43"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'166.32"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'9.28"
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("E28").Value = "  +4.37%  "
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("E33").Value = "  +2.02%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'2.05"
$ws.Range("E34").Value = "  +10.07%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0610"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("E37").Value = "  +9.70%  "
$ws.Range("E38").Value = "  +5.10%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").Value = "1.485.88"
$ws.Range("E40").Value = "  +0.99%  "
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("D42").Value = "'96.99"
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("D43").Value = "'2.85"
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("D44").Value = "'16.55"
$ws.Range("D45").Value = "'0.0923"
$ws.Range("E45").Value = "  +1.23%  "
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("D47").Value = "'4.13"
$ws.Range("E47").Value = "  +12.68%  "
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("E50").Value = "  -2.79%  "
$ws.Range("D51").Value = "2.242.40"
$ws.Range("E51").Value = "  +1.41%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
